# Updates the cryptos list (Coin/Link/Price/Volume) per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free inline approach: for columns that can hold numeric-looking
# text (Price column D), force a Text number format before assigning so
# Excel does not auto-convert the literal string into a number, then
# restore the default "Normal" style so formatting matches the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.338.22'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.26%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.458.22'
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.55'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.99%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.458.17'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.89%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.590'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +10.64%  '
$ws.Range("E10").Value = '  -2.10%  '
$ws.Range("E11").Value = '  +4.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.447'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.057.95'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.01%  '
$ws.Range("E14").Value = '  -2.66%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000195'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.71%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.97'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +7.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.379.36'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.460.66'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("E20").Value = '  +3.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '387.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.548'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.42'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.54%  '
$ws.Range("E25").Value = '  +0.64%  '
$ws.Range("E26").Value = '  +19.97%  '
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.180'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("E29").Value = '  +0.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.16'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.09%  '
$ws.Range("E31").Value = '  +8.63%  '
$ws.Range("E32").Value = '  -0.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.72'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.12'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.50'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '161.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0777'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.66%  '
$ws.Range("E40").Value = '  +0.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '27.58'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.930.87'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0321'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.98%  '
$ws.Range("E44").Value = '  +5.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.62'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.772'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.82'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.69%  '
$ws.Range("E48").Value = '  +2.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.23'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +17.21%  '
$ws.Range("E50").Value = '  +4.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.866'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.36%  '
